$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# The "Doing:" column (E) on the board only keeps its row-4 entry now;
# the three task cards that used to sit in E5:E7 were moved out.
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()

# Those cards were re-filed into the "Done:" column (G), continuing on
# from the existing G4 entry, shifting the rest of the Done list down.
$ws.Range("G7").Value = "Organizar o repositório no Github (José Pereira)"
$ws.Range("G8").Value = "Reunir o grupo para organizar a entrega da 2ª fase do projeto (José Pereira)"
$ws.Range("G9").Value = "Identificar 3 Design Patterns (Todos têm de fazer)"

# G8 used to be an empty, specially (unused-font) styled placeholder cell;
# now that it holds real content again it goes back to the default style.
$ws.Range("G8").Style = "Normal"

# View state: the sheet was re-zoomed and a different cell selected.
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
[void]$ws.Range("D13").Select()

# Best-effort: reflect the resized/repositioned workbook window too.
$excel.ActiveWindow.WindowState = -4143
$excel.ActiveWindow.Left = 22932
$excel.ActiveWindow.Top = -108
$excel.ActiveWindow.Width = 23256
$excel.ActiveWindow.Height = 12576
